$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.107.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6926"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07679"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3032"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.861.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7245"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.202"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.106.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.17%  "

$ws.Range("E17").Value = "  -4.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007767"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.02%  "

$ws.Range("E21").Value = "  +0.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.097.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.593"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.974"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1430"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.979"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.488"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.487"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.017"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05223"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.184"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.75%  "

$ws.Range("E36").Value = "  +2.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7003"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.65%  "

$ws.Range("E38").Value = "  -1.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01852"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.678"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9158"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.53%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.087.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.20%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.983"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4256"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.771"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.993.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.130"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.972"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.22%  "
